$wb = $excel.ActiveWorkbook

# Duplicate the "gaming_server" sheet (it already carries the identical layout,
# styles and merged header used by every "Test Case Summary" report sheet) and
# drop the copy right after "redeem_voucher" so the new "get_support" report
# lands in the same position as in the target workbook.
$template = $wb.Worksheets.Item("gaming_server")
$after = $wb.Worksheets.Item("redeem_voucher")
$template.Copy($null, $after)

$ws = $wb.Worksheets.Item($after.Index + 1)
$ws.Name = "get_support"

# Fill in the new report's own summary numbers.
$ws.Range("D2").Value = " Test Case Summary (06-03-24)"
$ws.Range("E3").Value = 20
$ws.Range("E4").Value = 19
$ws.Range("E5").Value = 0
$ws.Range("E6").Value = 1

$ws.Range("E6").Select()
